# Update cryptos list figures (Price / Volume(1h)) to match the latest
# scrape, per commit "Updated cryptos list on Sun Apr 30 07:59:17 UTC 2023
# with GitHub Actions".
#
# Column D ("Price") holds plain-text values (e.g. "326.02", "29.599.49").
# Several of the new prices (e.g. "326.02", "1.012") are syntactically
# valid numbers, and Excel's COM layer auto-converts a bare Value
# assignment like that into a numeric cell. To keep these cells as text
# (matching the original inlineStr/shared-string representation with no
# style change) we temporarily force the cell's number format to Text
# ("@") before assigning the value, then restore the cell style to
# "Normal" afterwards so no stray style index is left on the cell.
#
# Column E ("Volume(1h)") values already contain "%" and padding spaces,
# so Excel keeps them as text without any special handling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText($ws, $cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}


Set-PriceText $ws "D2" "29.599.49"
$ws.Range("E2").Value = "  +0.20%  "

Set-PriceText $ws "D3" "1.924.90"
$ws.Range("E3").Value = "  +0.49%  "

$ws.Range("E4").Value = "  +0.54%  "

Set-PriceText $ws "D5" "326.02"
$ws.Range("E5").Value = "  +0.03%  "

Set-PriceText $ws "D6" "1.012"
$ws.Range("E6").Value = "  +0.53%  "

Set-PriceText $ws "D7" "0.4818"
$ws.Range("E7").Value = "  -0.14%  "

Set-PriceText $ws "D8" "0.4063"
$ws.Range("E8").Value = "  -0.26%  "

Set-PriceText $ws "D9" "0.08228"
$ws.Range("E9").Value = "  +1.01%  "

Set-PriceText $ws "D11" "23.71"
$ws.Range("E11").Value = "  +1.20%  "

Set-PriceText $ws "D12" "1.927.96"
$ws.Range("E12").Value = "  -0.04%  "

Set-PriceText $ws "D13" "6.079"
$ws.Range("E13").Value = "  +1.35%  "

$ws.Range("E14").Value = "  +1.68%  "

Set-PriceText $ws "D15" "91.70"
$ws.Range("E15").Value = "  +1.62%  "

Set-PriceText $ws "D16" "0.06866"
$ws.Range("E16").Value = "  +1.18%  "

Set-PriceText $ws "D17" "1.013"
$ws.Range("E17").Value = "  +0.48%  "

$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("E20").Value = "  +0.46%  "

Set-PriceText $ws "D21" "29.592.60"
$ws.Range("E21").Value = "  +0.16%  "

Set-PriceText $ws "D22" "5.681"
$ws.Range("E22").Value = "  +0.93%  "

Set-PriceText $ws "D23" "11.97"
$ws.Range("E23").Value = "  +1.58%  "

Set-PriceText $ws "D24" "2.185"
$ws.Range("E24").Value = "  +0.08%  "

Set-PriceText $ws "D25" "2.155.52"
$ws.Range("E25").Value = "  -0.19%  "

Set-PriceText $ws "D26" "156.18"
$ws.Range("E26").Value = "  +0.42%  "

Set-PriceText $ws "D27" "6.468"
$ws.Range("E27").Value = "  +0.76%  "

Set-PriceText $ws "D28" "19.97"
$ws.Range("E28").Value = "  -0.24%  "

Set-PriceText $ws "D29" "2.093"
$ws.Range("E29").Value = "  -0.28%  "

Set-PriceText $ws "D30" "120.57"
$ws.Range("E30").Value = "  +0.77%  "

Set-PriceText $ws "D31" "1.015"
$ws.Range("E31").Value = "  -1.79%  "

Set-PriceText $ws "D32" "0.09629"
$ws.Range("E32").Value = "  +0.70%  "

Set-PriceText $ws "D33" "5.618"
$ws.Range("E33").Value = "  +2.01%  "

Set-PriceText $ws "D34" "3.583"
$ws.Range("E34").Value = "  +0.50%  "

$ws.Range("E35").Value = "  -1.12%  "

Set-PriceText $ws "D36" "0.06343"
$ws.Range("E36").Value = "  +4.10%  "

Set-PriceText $ws "D37" "0.02291"
$ws.Range("E37").Value = "  +1.00%  "

Set-PriceText $ws "D38" "1.188"
$ws.Range("E38").Value = "  +0.99%  "

Set-PriceText $ws "D39" "0.5957"
$ws.Range("E39").Value = "  +0.34%  "

Set-PriceText $ws "D40" "10.76"
$ws.Range("E40").Value = "  +0.33%  "

Set-PriceText $ws "D41" "7.878"
$ws.Range("E41").Value = "  -1.15%  "

$ws.Range("E42").Value = "  -0.14%  "

$ws.Range("E43").Value = "  -0.18%  "

Set-PriceText $ws "D44" "1.290"
$ws.Range("E44").Value = "  +0.57%  "

Set-PriceText $ws "D45" "12.45"
$ws.Range("E45").Value = "  -0.28%  "

Set-PriceText $ws "D46" "0.07499"
$ws.Range("E46").Value = "  -2.37%  "

Set-PriceText $ws "D47" "0.5555"
$ws.Range("E47").Value = "  -0.49%  "

$ws.Range("E48").Value = "  +0.45%  "

Set-PriceText $ws "D49" "119.44"
$ws.Range("E49").Value = "  +3.32%  "

Set-PriceText $ws "D50" "2.437"
$ws.Range("E50").Value = "  +3.45%  "

Set-PriceText $ws "D51" "72.21"
$ws.Range("E51").Value = "  -0.53%  "

Write-Output "cryptos list updated"
